# Updated cryptos list on Thu Feb 15 10:58:21 UTC 2024 with GitHub Actions
#
# This applies the per-cell Price (D) / Volume(1h) (E) refresh, plus the
# Monero / EnergySwap row swap (rows 42-43 trade places content-wise),
# against the already-open workbook ($excel.ActiveWorkbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text happens to look like a plain decimal number ---
# Excel normally auto-converts a bare numeric string assigned via .Value
# into a real Double, which would silently drop meaningful trailing zeros
# (e.g. "35.10" -> 35.1) and shift these cells from text to numbers. The
# source sheet stores every Price/Volume cell as literal text, so force the
# cell to Text format for the write, then put the style back to Normal so
# no stray number-format is left behind on the cell.
$numericTextCells = @(
    @{ Cell = "D5"; Value = "347.16" }
    @{ Cell = "D6"; Value = "116.68" }
    @{ Cell = "D9"; Value = "0.587" }
    @{ Cell = "D10"; Value = "43.27" }
    @{ Cell = "D14"; Value = "7.84" }
    @{ Cell = "D17"; Value = "0.894" }
    @{ Cell = "D19"; Value = "3.24" }
    @{ Cell = "D21"; Value = "13.49" }
    @{ Cell = "D23"; Value = "70.25" }
    @{ Cell = "D24"; Value = "270.24" }
    @{ Cell = "D26"; Value = "26.72" }
    @{ Cell = "D28"; Value = "10.27" }
    @{ Cell = "D31"; Value = "35.10" }
    @{ Cell = "D34"; Value = "0.0826" }
    @{ Cell = "D35"; Value = "0.0411" }
    @{ Cell = "D38"; Value = "18.95" }
    @{ Cell = "D40"; Value = "3.22" }
    @{ Cell = "D41"; Value = "2.71" }
    @{ Cell = "D42"; Value = "23.66" }
    @{ Cell = "D43"; Value = "128.39" }
    @{ Cell = "D49"; Value = "0.990" }
)
foreach ($item in $numericTextCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
    $rng.Style = "Normal"
}

# --- Remaining cells (coin names, links, multi-dot prices, percentages) ---
# These never look like a single Excel-parsable number, so a direct .Value
# assignment keeps them as text exactly as in the source file.
$textCells = @(
    @{ Cell = "D2"; Value = "52.295.71" }
    @{ Cell = "E2"; Value = "  +1.95%  " }
    @{ Cell = "D3"; Value = "2.799.07" }
    @{ Cell = "E3"; Value = "  +1.59%  " }
    @{ Cell = "E4"; Value = "  +0.01%  " }
    @{ Cell = "E5"; Value = "  +4.92%  " }
    @{ Cell = "E6"; Value = "  +1.09%  " }
    @{ Cell = "E8"; Value = "  -0.07%  " }
    @{ Cell = "E9"; Value = "  +2.26%  " }
    @{ Cell = "E10"; Value = "  +4.19%  " }
    @{ Cell = "E11"; Value = "  +3.30%  " }
    @{ Cell = "E12"; Value = "  -0.84%  " }
    @{ Cell = "E13"; Value = "  +1.61%  " }
    @{ Cell = "E14"; Value = "  +2.40%  " }
    @{ Cell = "D15"; Value = "3.237.59" }
    @{ Cell = "E15"; Value = "  +1.58%  " }
    @{ Cell = "D16"; Value = "2.785.70" }
    @{ Cell = "E16"; Value = "  +1.11%  " }
    @{ Cell = "E17"; Value = "  +0.49%  " }
    @{ Cell = "D18"; Value = "52.211.17" }
    @{ Cell = "E18"; Value = "  +1.81%  " }
    @{ Cell = "E19"; Value = "  +6.80%  " }
    @{ Cell = "E20"; Value = "  +3.83%  " }
    @{ Cell = "E21"; Value = "  -1.02%  " }
    @{ Cell = "D22"; Value = "0.0₃0982" }
    @{ Cell = "E22"; Value = "  +1.75%  " }
    @{ Cell = "E23"; Value = "  -0.08%  " }
    @{ Cell = "E24"; Value = "  -3.53%  " }
    @{ Cell = "E25"; Value = "  +5.93%  " }
    @{ Cell = "E26"; Value = "  -1.37%  " }
    @{ Cell = "E27"; Value = "  -0.04%  " }
    @{ Cell = "E28"; Value = "  -0.84%  " }
    @{ Cell = "E29"; Value = "  +0.74%  " }
    @{ Cell = "E30"; Value = "  -0.63%  " }
    @{ Cell = "E31"; Value = "  -1.79%  " }
    @{ Cell = "E32"; Value = "  +0.45%  " }
    @{ Cell = "E33"; Value = "  +1.54%  " }
    @{ Cell = "E34"; Value = "  -0.03%  " }
    @{ Cell = "E35"; Value = "  +16.07%  " }
    @{ Cell = "E36"; Value = "  +0.48%  " }
    @{ Cell = "E37"; Value = "  -0.04%  " }
    @{ Cell = "E38"; Value = "  -2.32%  " }
    @{ Cell = "E39"; Value = "  -1.02%  " }
    @{ Cell = "E40"; Value = "  -0.43%  " }
    @{ Cell = "E41"; Value = "  +22.08%  " }
    @{ Cell = "B42"; Value = "EnergySwap" }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" }
    @{ Cell = "E42"; Value = "  -0.41%  " }
    @{ Cell = "B43"; Value = "Monero" }
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" }
    @{ Cell = "E43"; Value = "  -0.87%  " }
    @{ Cell = "E44"; Value = "  +1.81%  " }
    @{ Cell = "E45"; Value = "  +0.38%  " }
    @{ Cell = "E46"; Value = "  -2.04%  " }
    @{ Cell = "E47"; Value = "  +4.93%  " }
    @{ Cell = "D48"; Value = "2.074.84" }
    @{ Cell = "E48"; Value = "  -2.08%  " }
    @{ Cell = "E49"; Value = "  +18.41%  " }
    @{ Cell = "E50"; Value = "  -0.21%  " }
    @{ Cell = "E51"; Value = "  -1.23%  " }
)
foreach ($item in $textCells) {
    $ws.Range($item.Cell).Value = $item.Value
}
